$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Kharagauli Municipality - Area" sheet used to show a 3-year comparison
# (1989 / 2002 / 2014). This edit trims it down to a single-year (2014)
# snapshot: columns C:D (1989, 2002) are dropped, and the remaining small
# table is re-laid-out/re-labelled to its final compact form.
# ---------------------------------------------------------------------------

$xlPasteFormats = -4122
$xlPasteValues  = -4163
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlMedium = -4138

# 1) Drop the 1989 and 2002 columns entirely - only the 2014 column (now B)
#    survives, carrying the (identical) 913.9 value forward.
$ws.Range("C1:D6").EntireColumn.Delete()

# 2) Stash the formatting of the soon-to-be-overwritten A5 ("year row" label
#    cell, currently blank) in a scratch cell so it isn't lost once A5 gets
#    the "Area" label copied onto it in step 4.
$ws.Range("A5").Copy()
$ws.Range("Z1").PasteSpecial($xlPasteFormats)

# 3) The unit caption "(sq. km)" (currently A6) moves up to become the
#    section header directly under the title (new A3).
$ws.Range("A6").Copy()
$ws.Range("A3").PasteSpecial($xlPasteFormats)
$ws.Range("A3").Value = "(sq. km)"

# 4) The row caption "Area" (currently A4) moves down to become the data
#    row's label (new A5).
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial($xlPasteFormats)
$ws.Range("A5").Value = "Area"

# 5) The blank label cell that used to sit above the data row (old A5) is
#    now the blank cell above the year header (new A4).
$ws.Range("Z1").Copy()
$ws.Range("A4").PasteSpecial($xlPasteFormats)
$ws.Range("A4").ClearContents()
$ws.Range("Z1").Clear()

# 6) Year header: old B5 (1989) becomes new B4, with its value corrected to
#    2014 (the only year that remains).
$ws.Range("B5").Copy()
$ws.Range("B4").PasteSpecial($xlPasteFormats)
$ws.Range("B4").Value = 2014

# 7) Data value: old B6 (913.9) becomes new B5. Its format starts from B6's
#    but - since column B is now the last column of the table - the right
#    edge must become the "closing" medium border instead of the old
#    interior thin/none edge.
$ws.Range("B6").Copy()
$ws.Range("B5").PasteSpecial($xlPasteFormats)
$ws.Range("B5").Value = 913.9
$ws.Range("B5").Borders.Item($xlEdgeRight).LineStyle = 1
$ws.Range("B5").Borders.Item($xlEdgeRight).Weight = $xlMedium

# 8) Clean up the cells whose content moved elsewhere / disappeared.
$ws.Range("B1").ClearContents()
$ws.Range("B1").ClearFormats()
$ws.Range("A2").ClearContents()
$ws.Range("A2").ClearFormats()
$ws.Range("B2").ClearContents()
$ws.Range("B2").ClearFormats()
$ws.Range("B3").ClearContents()
$ws.Range("B3").ClearFormats()
$ws.Range("A6").ClearContents()
$ws.Range("A6").ClearFormats()
$ws.Range("B6").ClearContents()
$ws.Range("B6").ClearFormats()

# 9) Row heights: every row in the compact table is a taller 20.1pt (was a
#    mix of 15 / 13.5 / default 12.75).
$ws.Rows.Item(1).RowHeight = 20.1
$ws.Rows.Item(2).RowHeight = 20.1
$ws.Rows.Item(3).RowHeight = 20.1
$ws.Rows.Item(4).RowHeight = 20.1
$ws.Rows.Item(5).RowHeight = 20.1
$ws.Rows.Item(6).RowHeight = 20.1

Write-Host "done"
